$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 360794.72
$ws.Range("J17").Value = 375801.88
$ws.Range("L17").Value = 1127405.64
$ws.Range("N17").Value = -1127741.64
# Row 33
$ws.Range("H33").Value = 236
$ws.Range("I33").Value = 190.6
$ws.Range("J33").Value = 690
$ws.Range("K33").Value = 190.6
$ws.Range("L33").Value = 690
$ws.Range("M33").Value = 38.40000000000001
$ws.Range("N33").Value = -1148
# Row 53
$ws.Range("H53").Value = 193.3077
$ws.Range("I53").Value = 93.666664
$ws.Range("J53").Value = 223.2
$ws.Range("K53").Value = 93.666664
$ws.Range("L53").Value = 223.2
$ws.Range("M53").Value = 543.333336
$ws.Range("N53").Value = -1497.2
# Row 74
$ws.Range("H74").Value = 3547.7222
$ws.Range("I74").Value = 3433.25
$ws.Range("K74").Value = 3433.25
$ws.Range("M74").Value = -2497.25
# Row 77
$ws.Range("H77").Value = 3547.7222
$ws.Range("I77").Value = 3433.25
$ws.Range("K77").Value = 17166.25
$ws.Range("M77").Value = -12486.25
# Row 137
$ws.Range("H137").Value = 2424.2334
$ws.Range("I137").Value = 2435.0952
$ws.Range("J137").Value = 2398.889
$ws.Range("K137").Value = 7305.285600000001
$ws.Range("L137").Value = 7196.667
$ws.Range("M137").Value = -4755.285600000001
$ws.Range("N137").Value = -12296.667

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 163.75
$ws.Range("I5").Value = 52.5
$ws.Range("J5").Value = 275
$ws.Range("K5").Value = 52.5
$ws.Range("L5").Value = 275
$ws.Range("M5").Value = 59.5
$ws.Range("N5").Value = -499
# Row 31
$ws.Range("H31").Value = 17587.6
$ws.Range("I31").Value = 4300
$ws.Range("K31").Value = 4300
$ws.Range("M31").Value = -4006
# Row 32
$ws.Range("H32").Value = 5128.3936
$ws.Range("I32").Value = 3924.3396
$ws.Range("K32").Value = 3924.3396
$ws.Range("M32").Value = -3637.3396
# Row 61
$ws.Range("H61").Value = 609.7907
$ws.Range("I61").Value = 596.425
$ws.Range("J61").Value = 788
$ws.Range("K61").Value = 596.425
$ws.Range("L61").Value = 788
$ws.Range("M61").Value = -384.425
$ws.Range("N61").Value = -1212
# Row 63
$ws.Range("H63").Value = 4163.3335
$ws.Range("J63").Value = 4100
$ws.Range("L63").Value = 4100
$ws.Range("N63").Value = -5472
# Row 66
$ws.Range("H66").Value = 4163.3335
$ws.Range("J66").Value = 4100
$ws.Range("L66").Value = 20500
$ws.Range("N66").Value = -27364
# Row 132
$ws.Range("H132").Value = 5634.6665
$ws.Range("I132").Value = 6295.0977
$ws.Range("J132").Value = 2926.9
$ws.Range("K132").Value = 18885.2931
$ws.Range("L132").Value = 8780.700000000001
$ws.Range("M132").Value = -16355.2931
$ws.Range("N132").Value = -13840.7
# Row 136
$ws.Range("H136").Value = 609.7907
$ws.Range("I136").Value = 596.425
$ws.Range("J136").Value = 788
$ws.Range("K136").Value = 1789.275
$ws.Range("L136").Value = 2364
$ws.Range("M136").Value = 760.7250000000001
$ws.Range("N136").Value = -7464
# Row 139
$ws.Range("H139").Value = 83111
$ws.Range("J139").Value = 96000
$ws.Range("L139").Value = 96000
$ws.Range("N139").Value = -106280

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 163.75
$ws.Range("I4").Value = 52.5
$ws.Range("J4").Value = 275
$ws.Range("K4").Value = 52.5
$ws.Range("L4").Value = 275
$ws.Range("M4").Value = 62.5
$ws.Range("N4").Value = -505
# Row 15
$ws.Range("H15").Value = 25005.834
$ws.Range("J15").Value = 25005.834
$ws.Range("L15").Value = 25005.834
$ws.Range("N15").Value = -25459.834
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 82
$ws.Range("H82").Value = 18028.3
$ws.Range("J82").Value = 33056.6
$ws.Range("L82").Value = 33056.6
$ws.Range("N82").Value = -33822.6
# Row 85
$ws.Range("H85").Value = 18028.3
$ws.Range("J85").Value = 33056.6
$ws.Range("L85").Value = 33056.6
$ws.Range("N85").Value = -35708.6

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1134.7142
$ws.Range("I58").Value = 1143.2051
$ws.Range("J58").Value = 1101.6
$ws.Range("K58").Value = 1143.2051
$ws.Range("L58").Value = 1101.6
$ws.Range("M58").Value = -940.2050999999999
$ws.Range("N58").Value = -1507.6
# Row 68
$ws.Range("H68").Value = 32250.1
$ws.Range("J68").Value = 32250.1
$ws.Range("L68").Value = 32250.1
$ws.Range("N68").Value = -33748.1
# Row 71
$ws.Range("H71").Value = 32250.1
$ws.Range("J71").Value = 32250.1
$ws.Range("L71").Value = 96750.29999999999
$ws.Range("N71").Value = -104238.3
# Row 74
$ws.Range("H74").Value = 35522.375
$ws.Range("J74").Value = 35522.375
$ws.Range("L74").Value = 35522.375
$ws.Range("N74").Value = -37270.375
# Row 77
$ws.Range("H77").Value = 35522.375
$ws.Range("J77").Value = 35522.375
$ws.Range("L77").Value = 106567.125
$ws.Range("N77").Value = -115303.125
# Row 134
$ws.Range("H134").Value = 1090.2894
$ws.Range("I134").Value = 924.44116
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 2773.32348
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -238.32348
$ws.Range("N134").Value = -12570
# Row 136
$ws.Range("H136").Value = 1134.7142
$ws.Range("I136").Value = 1143.2051
$ws.Range("J136").Value = 1101.6
$ws.Range("K136").Value = 3429.615299999999
$ws.Range("L136").Value = 3304.8
$ws.Range("M136").Value = -879.6152999999995
$ws.Range("N136").Value = -8404.799999999999
# Row 138
$ws.Range("H138").Value = 46546.668
$ws.Range("J138").Value = 46546.668
$ws.Range("L138").Value = 46546.668
$ws.Range("N138").Value = -56826.668

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1789672
$ws.Range("I131").Value = 386.66666
$ws.Range("J131").Value = 2263306.2
$ws.Range("K131").Value = 1159.99998
$ws.Range("L131").Value = 6789918.600000001
$ws.Range("M131").Value = 3880.00002
$ws.Range("N131").Value = -6799998.600000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5097.393
$ws.Range("I70").Value = 5063.4116
$ws.Range("K70").Value = 5063.4116
$ws.Range("M70").Value = -4793.4116
# Row 73
$ws.Range("H73").Value = 5097.393
$ws.Range("I73").Value = 5063.4116
$ws.Range("K73").Value = 5063.4116
$ws.Range("M73").Value = -4127.4116
# Row 80
$ws.Range("H80").Value = 2444.4443
$ws.Range("I80").Value = 2360
$ws.Range("J80").Value = 2550
$ws.Range("K80").Value = 2360
$ws.Range("L80").Value = 2550
$ws.Range("M80").Value = -1362
$ws.Range("N80").Value = -4546
# Row 83
$ws.Range("H83").Value = 2444.4443
$ws.Range("I83").Value = 2360
$ws.Range("J83").Value = 2550
$ws.Range("K83").Value = 11800
$ws.Range("L83").Value = 12750
$ws.Range("M83").Value = -6808
$ws.Range("N83").Value = -22734
# Row 113
$ws.Range("H113").Value = 905.26666
$ws.Range("I113").Value = 850.2
$ws.Range("J113").Value = 1015.4
$ws.Range("K113").Value = 850.2
$ws.Range("L113").Value = 1015.4
$ws.Range("M113").Value = 1319.8
$ws.Range("N113").Value = -5355.4
# Row 122
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -13900
# Row 132
$ws.Range("H132").Value = 8469.666999999999
$ws.Range("I132").Value = 11285.429
$ws.Range("J132").Value = 1899.5555
$ws.Range("K132").Value = 33856.287
$ws.Range("L132").Value = 5698.666499999999
$ws.Range("M132").Value = -31326.287
$ws.Range("N132").Value = -10758.6665

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1570.5714
$ws.Range("I82").Value = 630.1
$ws.Range("J82").Value = 2425.5454
$ws.Range("K82").Value = 630.1
$ws.Range("L82").Value = 2425.5454
$ws.Range("M82").Value = -269.1
$ws.Range("N82").Value = -3147.5454
# Row 85
$ws.Range("H85").Value = 1570.5714
$ws.Range("I85").Value = 630.1
$ws.Range("J85").Value = 2425.5454
$ws.Range("K85").Value = 630.1
$ws.Range("L85").Value = 2425.5454
$ws.Range("M85").Value = 617.9
$ws.Range("N85").Value = -4921.5454
# Row 132
$ws.Range("H132").Value = 5524.5
$ws.Range("I132").Value = 5924.409
$ws.Range("J132").Value = 4424.75
$ws.Range("K132").Value = 17773.227
$ws.Range("L132").Value = 13274.25
$ws.Range("M132").Value = -15243.227
$ws.Range("N132").Value = -18334.25
